$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Session 11 (row 12): add Vorbereitung (prep) material
$ws.Range("D12").Value = "prep/p11.html"

# Session 10 (row 11): add Abgaben (exercise) entry and Folien (slides) link
$ws.Range("F11").Value = "exercises/e10.html"
$ws.Range("E11").Value = "slides/slides.html#/sitzung-10-human-agency-wie-menschen-sich-auf-algorithmisch-kuratierten-plattformen-unterhalten"

# Update the selected cell to match the new state
$ws.Range("E11").Select()
